# Updates cryptos list: prices and 1h volume % changes for each coin row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.916.88"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.876.01"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -3.79%  "
$ws.Range("D6").Value = "'242.62"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.3148"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "'0.07232"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").Value = "'24.66"
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("D11").Value = "'0.08341"
$ws.Range("E11").Value = "  -2.35%  "
$ws.Range("D12").Value = "'0.7499"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.882.62"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.392"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "'92.37"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "29.934.19"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "'6.103"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "'247.16"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "'13.57"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "'0.000007843"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "2.141.27"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").Value = "'8.011"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'0.1558"
$ws.Range("E25").Value = "  -5.31%  "
$ws.Range("D26").Value = "'9.292"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").Value = "'165.20"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").Value = "'18.67"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").Value = "'2.025"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").Value = "'1.495"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("D31").Value = "'4.612"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("D32").Value = "'1.537"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("E33").Value = "  +3.93%  "
$ws.Range("D34").Value = "'0.05323"
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("D36").Value = "'0.7502"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").Value = "'0.9988"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "'2.699"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "'2.755"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").Value = "'0.4529"
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.109.37"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'6.129"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "'72.39"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").Value = "'0.8646"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").Value = "'104.53"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'7.590"
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").Value = "'9.530"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").Value = "2.038.47"
$ws.Range("E51").Value = "  -0.39%  "
